$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 813, shifting existing rows 813:840 down to 814:841
$ws.Rows.Item(813).Insert()

# Populate the newly inserted row 813 with the new record
$ws.Cells.Item(813, 1).Value = 8
$ws.Cells.Item(813, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(813, 3).Value = "Coquimbo"
$ws.Cells.Item(813, 4).Value = Get-Date -Year 2023 -Month 5 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(813, 5).Value = 4
$ws.Cells.Item(813, 6).Value = 100112024
$ws.Cells.Item(813, 7).Value = "Choclo"
$ws.Cells.Item(813, 8).Value = "Dulce o Americano"
$ws.Cells.Item(813, 9).Value = "Primera"
$ws.Cells.Item(813, 10).Value = 12000
$ws.Cells.Item(813, 11).Value = 230
$ws.Cells.Item(813, 12).Value = 250
$ws.Cells.Item(813, 13).Value = 240
$ws.Cells.Item(813, 14).Value = "`$/unidad"
$ws.Cells.Item(813, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(813, 16).Value = 240
$ws.Cells.Item(813, 17).Value = 1
$ws.Cells.Item(813, 18).Value = "Hortaliza"
